$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.216.44"
$ws.Range("E2").Value = "  +0.08%  "
$ws.Range("D3").Value = "1.851.53"
$ws.Range("E3").Value = "  -0.33%  "
$c = $ws.Range("D4")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.Style = $origStyle
$ws.Range("E4").Value = "  +0.15%  "
$c = $ws.Range("D5")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.6961"
$c.Style = $origStyle
$ws.Range("E5").Value = "  +0.71%  "
$ws.Range("E6").Value = "  -0.37%  "
$ws.Range("E7").Value = "  +0.15%  "
$c = $ws.Range("D8")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.07893"
$c.Style = $origStyle
$ws.Range("E8").Value = "  +1.94%  "
$c = $ws.Range("D9")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.3010"
$c.Style = $origStyle
$ws.Range("E9").Value = "  -1.57%  "
$c = $ws.Range("D10")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "23.48"
$c.Style = $origStyle
$ws.Range("E10").Value = "  +0.74%  "
$c = $ws.Range("D11")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.08149"
$c.Style = $origStyle
$ws.Range("E11").Value = "  +0.97%  "
$ws.Range("D12").Value = "1.851.08"
$ws.Range("E12").Value = "  -0.61%  "
$c = $ws.Range("D13")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "5.170"
$c.Style = $origStyle
$ws.Range("E13").Value = "  -0.70%  "
$ws.Range("E14").Value = "  -2.91%  "
$c = $ws.Range("D15")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "89.35"
$c.Style = $origStyle
$ws.Range("E15").Value = "  -0.24%  "
$ws.Range("D16").Value = "29.244.71"
$ws.Range("E16").Value = "  +0.16%  "
$c = $ws.Range("D17")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "5.791"
$c.Style = $origStyle
$ws.Range("E17").Value = "  +0.70%  "
$c = $ws.Range("D18")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.000007797"
$c.Style = $origStyle
$c = $ws.Range("D19")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "13.15"
$c.Style = $origStyle
$ws.Range("E19").Value = "  -1.02%  "
$c = $ws.Range("D20")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "235.73"
$c.Style = $origStyle
$ws.Range("E20").Value = "  +0.02%  "
$ws.Range("E21").Value = "  +0.20%  "
$ws.Range("D22").Value = "2.101.71"
$ws.Range("E22").Value = "  -0.06%  "
$c = $ws.Range("D23")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.Style = $origStyle
$ws.Range("E23").Value = "  +0.13%  "
$c = $ws.Range("D24")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "7.513"
$c.Style = $origStyle
$ws.Range("E24").Value = "  +0.71%  "
$c = $ws.Range("D25")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "162.44"
$c.Style = $origStyle
$ws.Range("E25").Value = "  +0.32%  "
$c = $ws.Range("D26")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "8.838"
$c.Style = $origStyle
$ws.Range("E26").Value = "  -1.67%  "
$c = $ws.Range("D27")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.1409"
$c.Style = $origStyle
$ws.Range("E27").Value = "  -2.36%  "
$c = $ws.Range("D28")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "18.01"
$c.Style = $origStyle
$c = $ws.Range("D29")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.909"
$c.Style = $origStyle
$ws.Range("E29").Value = "  -2.74%  "
$c = $ws.Range("D30")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.406"
$c.Style = $origStyle
$ws.Range("E30").Value = "  +0.41%  "
$c = $ws.Range("D31")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.469"
$c.Style = $origStyle
$ws.Range("E31").Value = "  -1.34%  "
$ws.Range("E32").Value = "  -4.73%  "
$c = $ws.Range("D33")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "4.010"
$c.Style = $origStyle
$ws.Range("E33").Value = "  -0.41%  "
$c = $ws.Range("D34")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.05140"
$c.Style = $origStyle
$ws.Range("E34").Value = "  -0.81%  "
$ws.Range("E35").Value = "  -2.17%  "
$c = $ws.Range("D36")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.7066"
$c.Style = $origStyle
$ws.Range("E36").Value = "  +0.16%  "
$c = $ws.Range("D37")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.9981"
$c.Style = $origStyle
$ws.Range("E37").Value = "  -2.54%  "
$c = $ws.Range("D38")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "2.681"
$c.Style = $origStyle
$ws.Range("E38").Value = "  +0.34%  "
$c = $ws.Range("D39")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.01840"
$c.Style = $origStyle
$ws.Range("E39").Value = "  -0.65%  "
$c = $ws.Range("D40")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "2.708"
$c.Style = $origStyle
$ws.Range("E40").Value = "  +0.95%  "
$ws.Range("D41").Value = "1.149.66"
$ws.Range("E41").Value = "  +4.81%  "
$c = $ws.Range("D42")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.9258"
$c.Style = $origStyle
$ws.Range("E42").Value = "  +0.16%  "
$c = $ws.Range("D43")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "5.965"
$c.Style = $origStyle
$ws.Range("E43").Value = "  -0.34%  "
$c = $ws.Range("D44")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.4225"
$c.Style = $origStyle
$ws.Range("E44").Value = "  -1.54%  "
$c = $ws.Range("D45")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "69.93"
$c.Style = $origStyle
$ws.Range("E45").Value = "  -0.99%  "
$ws.Range("E46").Value = "  +0.13%  "
$c = $ws.Range("D47")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "102.46"
$c.Style = $origStyle
$ws.Range("E47").Value = "  +0.30%  "
$c = $ws.Range("D48")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.5294"
$c.Style = $origStyle
$ws.Range("E48").Value = "  -2.83%  "
$c = $ws.Range("D49")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.731"
$c.Style = $origStyle
$ws.Range("E49").Value = "  -3.55%  "
$c = $ws.Range("D50")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "9.102"
$c.Style = $origStyle
$c = $ws.Range("D51")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "6.932"
$c.Style = $origStyle
$ws.Range("E51").Value = "  -1.15%  "
